$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C29").Value = 59
$ws.Range("D29").Value = 12
$ws.Range("E29").Value = 47
$ws.Range("F29").Value = 2.065404475043029
